# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Replaces the "VALERIA FRANCO HERNANDEZ" worker block with a new
# "LUIS FELIPE BRU TABORDA" block (periods 2001-2002) inserted before the
# existing "YESSI ARNALDO MARTINEZ IRIARTE" block (periods 2201-2210), and
# updates the summary figures at the top of the sheet accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Preserve the special "last row of the table" border formatting that
#    currently lives on row 29 by copying it onto row 27 (the row that
#    will become the new last data row once rows 28:29 are removed).
# ---------------------------------------------------------------------
$ws.Range("B29:J29").Copy()
$ws.Range("B27:J27").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 2) Remove the two rows belonging to "VALERIA FRANCO HERNANDEZ"
#    (periods 2507 / 2505). This also shifts the signature-block rows
#    (34/35 -> 32/33) up automatically, fixing dimension/mergeCells.
# ---------------------------------------------------------------------
$ws.Rows("28:29").Delete()

# ---------------------------------------------------------------------
# 3) Update the summary header values.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 412567    # VALOR MORA
$ws.Range("C13").Value = 2         # Cant. Trabajadores
$ws.Range("F13").Value = 12        # Cant. Periodos

# ---------------------------------------------------------------------
# 4) Rewrite the worker detail rows (16-27):
#    rows 16-17 -> LUIS FELIPE BRU TABORDA, periods 2001 / 2002
#    rows 18-27 -> YESSI ARNALDO MARTINEZ IRIARTE, periods 2201-2210
# ---------------------------------------------------------------------

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047495227"
$ws.Range("D16").Value = "LUIS FELIPE BRU TABORDA"
$ws.Range("E16").Value = "2001"
$ws.Range("F16").Value = 14045
$ws.Range("G16").Value = 877803

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047495227"
$ws.Range("D17").Value = "LUIS FELIPE BRU TABORDA"
$ws.Range("E17").Value = "2002"
$ws.Range("F17").Value = 35112
$ws.Range("G17").Value = 877803

$periods = @("2201","2202","2203","2204","2205","2206","2207","2208","2209","2210")
$row = 18
foreach ($p in $periods) {
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "1047430916"
    $ws.Range("D$row").Value = "YESSI ARNALDO MARTINEZ IRIARTE"
    $ws.Range("E$row").Value = $p
    $ws.Range("F$row").Value = 36341
    $ws.Range("G$row").Value = 908526
    $row = $row + 1
}
